# Apply the SimParameters treatment-effect-heterogeneity edit:
# - Moderate/High severity RR for Abortion moves from B9/B10 into new C9/C10
#   (the "ratio scale" multiplier), and B9/B10 are repurposed to hold the
#   plain severity level (2/3) used elsewhere in the workbook.
# - Downstream sheets (potential_preg_untrt, potential_preg_trt) recompute
#   automatically because their cells formula-reference SimParameters!B9/B10.
# - Finally, make SimParameters the active sheet/tab with B11 selected
#   (matching the author's end-of-edit UI state), since it had been
#   potential_preg_untrt/C9 before.

$wb = $excel.ActiveWorkbook
$simParams = $wb.Worksheets.Item("SimParameters")

# Preserve the old ratio-scale RR values by moving them into column C,
# then overwrite column B with the new severity-level values.
$simParams.Range("C9").Value2 = 1.25
$simParams.Range("B9").Value2 = 2

$simParams.Range("C10").Value2 = 1.5
$simParams.Range("B10").Value2 = 3

# Recalculate so every dependent formula cell picks up the new values.
$excel.Calculate()

# Update the active sheet / selection to match the saved UI state.
$simParams.Activate()
$simParams.Range("B11").Select()
